$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Cells.Item(1,1).Value = "Datos actualizados a 2 de Agosto de 2020 a las 10:47"

# Helper to set a full data row (columns B..H) at once
function Set-RowData($Row, $B, $C, $D, $E, $F, $G, $H) {
    $ws.Cells.Item($Row,2).Value = $B
    $ws.Cells.Item($Row,3).Value = $C
    $ws.Cells.Item($Row,4).Value = $D
    $ws.Cells.Item($Row,5).Value = $E
    $ws.Cells.Item($Row,6).Value = $F
    $ws.Cells.Item($Row,7).Value = $G
    $ws.Cells.Item($Row,8).Value = $H
}

# Row 4  - Estados Unidos
Set-RowData 4   4764588 270  2363165 2243518 0 7  157905
# Row 6  - India
Set-RowData 6   1756220 4301 1148161 570623  0 33 37436
# Row 7  - Rusia
Set-RowData 7   850870  5427 650173  186569  0 70 14128
# Row 28 - Filipinas
Set-RowData 28  103185  5032 65557   35569   0 20 2059
# Row 36 - Israel
Set-RowData 36  72315   97   45631   26153   0 5  531
# Row 45 - Singapur
Set-RowData 45  52825   313  46740   6058    0 0  27
# Row 63 - Moldavia
Set-RowData 63  25113   0    17816   6508    0 1  789
# Row 73 - now El Salvador (updated figures, list re-sorted by total cases)
Set-RowData 73  17448   398  8634    8347    0 8  467
# Row 74 - now Camerun (previous Camerun figures, list re-sorted by total cases)
Set-RowData 74  17255   0    15320   1544    0 0  391
# Row 88 - Malasia
Set-RowData 88  8999    14   8664    210     0 0  125
# Row 126 - Eslovenia
Set-RowData 126 2180    9    1826    234     0 1  120
# Row 127 - Lituania
Set-RowData 127 2110    17   1644    386     0 0  80

# Row 73/74 country names swap in-place: the list is sorted descending by total
# cases (column B). El Salvador's total overtook Camerun's, so their rows swap
# identity while their sheet rows (73/74) remain fixed.
$ws.Cells.Item(73,1).Value = "El Salvador"
$ws.Cells.Item(74,1).Value = "Camerun"
